# Auto-generated edit script applying numeric corrections to Leve profit tables
# across the ALC, ARM, BSM, CRP, CUL, GSM, and LTW worksheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 70024
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H35").Value = 70024
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H40").Value = 994
$ws.Range("I40").Value = 1037.5
$ws.Range("J40").Value = 955.3333
$ws.Range("K40").Value = 1037.5
$ws.Range("L40").Value = 955.3333
$ws.Range("M40").Value = -862.5
$ws.Range("N40").Value = -1305.3333
$ws.Range("H132").Value = 2633152.8
$ws.Range("I132").Value = 2942651.2
$ws.Range("J132").Value = 2415
$ws.Range("K132").Value = 8827953.600000001
$ws.Range("L132").Value = 7245
$ws.Range("M132").Value = -8825423.600000001
$ws.Range("N132").Value = -12305
$ws.Range("H137").Value = 2050.8853
$ws.Range("I137").Value = 2053.347
$ws.Range("K137").Value = 6160.041000000001
$ws.Range("M137").Value = -3610.041000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4706.39
$ws.Range("I32").Value = 4418.506
$ws.Range("J32").Value = 6111.9414
$ws.Range("K32").Value = 4418.506
$ws.Range("L32").Value = 6111.9414
$ws.Range("M32").Value = -4131.506
$ws.Range("N32").Value = -6685.9414
$ws.Range("H61").Value = 3241.7188
$ws.Range("I61").Value = 1187.3077
$ws.Range("K61").Value = 1187.3077
$ws.Range("M61").Value = -975.3077000000001
$ws.Range("H74").Value = 990.0454999999999
$ws.Range("I74").Value = 944.63635
$ws.Range("J74").Value = 1035.4546
$ws.Range("K74").Value = 944.63635
$ws.Range("L74").Value = 1035.4546
$ws.Range("M74").Value = -70.63634999999999
$ws.Range("N74").Value = -2783.4546
$ws.Range("H77").Value = 990.0454999999999
$ws.Range("I77").Value = 944.63635
$ws.Range("J77").Value = 1035.4546
$ws.Range("K77").Value = 4723.18175
$ws.Range("L77").Value = 5177.273
$ws.Range("M77").Value = -355.1817499999997
$ws.Range("N77").Value = -13913.273
$ws.Range("H136").Value = 3241.7188
$ws.Range("I136").Value = 1187.3077
$ws.Range("K136").Value = 3561.9231
$ws.Range("M136").Value = -1011.9231
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3694.182
$ws.Range("I107").Value = 2177.875
$ws.Range("J107").Value = 7737.6665
$ws.Range("K107").Value = 2177.875
$ws.Range("L107").Value = 7737.6665
$ws.Range("M107").Value = -257.875
$ws.Range("N107").Value = -11577.6665
$ws.Range("H134").Value = 1410.4032
$ws.Range("I134").Value = 1070.9824
$ws.Range("J134").Value = 5279.8
$ws.Range("K134").Value = 3212.947200000001
$ws.Range("L134").Value = 15839.4
$ws.Range("M134").Value = -677.9472000000005
$ws.Range("N134").Value = -20909.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2936.7556
$ws.Range("I31").Value = 1750.7273
$ws.Range("J31").Value = 4071.2173
$ws.Range("K31").Value = 1750.7273
$ws.Range("L31").Value = 4071.2173
$ws.Range("M31").Value = -1455.7273
$ws.Range("N31").Value = -4661.2173
$ws.Range("H34").Value = 2936.7556
$ws.Range("I34").Value = 1750.7273
$ws.Range("J34").Value = 4071.2173
$ws.Range("K34").Value = 1750.7273
$ws.Range("L34").Value = 4071.2173
$ws.Range("M34").Value = -1548.7273
$ws.Range("N34").Value = -4475.2173
$ws.Range("H39").Value = 25583.5
$ws.Range("I39").Value = 4990
$ws.Range("K39").Value = 4990
$ws.Range("M39").Value = -4599
$ws.Range("H49").Value = 25583.5
$ws.Range("I49").Value = 4990
$ws.Range("K49").Value = 4990
$ws.Range("M49").Value = -4808
$ws.Range("H58").Value = 9806387
$ws.Range("I58").Value = 1492.9143
$ws.Range("K58").Value = 1492.9143
$ws.Range("M58").Value = -1289.9143
$ws.Range("H122").Value = 3641.25
$ws.Range("I122").Value = 2846.125
$ws.Range("J122").Value = 4436.375
$ws.Range("K122").Value = 8538.375
$ws.Range("L122").Value = 13309.125
$ws.Range("M122").Value = -6088.375
$ws.Range("N122").Value = -18209.125
$ws.Range("H134").Value = 1360.0238
$ws.Range("I134").Value = 614.1389
$ws.Range("J134").Value = 5835.3335
$ws.Range("K134").Value = 1842.4167
$ws.Range("L134").Value = 17506.0005
$ws.Range("M134").Value = 692.5832999999998
$ws.Range("N134").Value = -22576.0005
$ws.Range("H135").Value = 26148.428
$ws.Range("J135").Value = 26148.428
$ws.Range("L135").Value = 26148.428
$ws.Range("N135").Value = -36288.428
$ws.Range("H136").Value = 9806387
$ws.Range("I136").Value = 1492.9143
$ws.Range("K136").Value = 4478.742899999999
$ws.Range("M136").Value = -1928.742899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3900
$ws.Range("J54").Value = 3900
$ws.Range("L54").Value = 11700
$ws.Range("N54").Value = -12818
$ws.Range("H64").Value = 2408.2144
$ws.Range("I64").Value = 1403
$ws.Range("J64").Value = 2966.6667
$ws.Range("K64").Value = 4209
$ws.Range("L64").Value = 8900.000100000001
$ws.Range("M64").Value = -3939
$ws.Range("N64").Value = -9440.000100000001
$ws.Range("H67").Value = 2408.2144
$ws.Range("I67").Value = 1403
$ws.Range("J67").Value = 2966.6667
$ws.Range("K67").Value = 4209
$ws.Range("L67").Value = 8900.000100000001
$ws.Range("M67").Value = -3273
$ws.Range("N67").Value = -10772.0001
$ws.Range("H114").Value = 900.8182
$ws.Range("I114").Value = 267.66666
$ws.Range("J114").Value = 1138.25
$ws.Range("K114").Value = 802.9999799999999
$ws.Range("L114").Value = 3414.75
$ws.Range("M114").Value = 2451.00002
$ws.Range("N114").Value = -9922.75
$ws.Range("H129").Value = 24081.738
$ws.Range("I129").Value = 2550
$ws.Range("J129").Value = 57575.555
$ws.Range("K129").Value = 7650
$ws.Range("L129").Value = 172726.665
$ws.Range("M129").Value = -2650
$ws.Range("N129").Value = -182726.665
$ws.Range("H131").Value = 1928.2759
$ws.Range("I131").Value = 5972.5
$ws.Range("J131").Value = 1281.2
$ws.Range("K131").Value = 17917.5
$ws.Range("L131").Value = 3843.6
$ws.Range("M131").Value = -12877.5
$ws.Range("N131").Value = -13923.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 51510
$ws.Range("J48").Value = 63353.332
$ws.Range("L48").Value = 63353.332
$ws.Range("N48").Value = -64323.332
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2979.5881
$ws.Range("I132").Value = 1939.7333
$ws.Range("J132").Value = 3800.5264
$ws.Range("K132").Value = 5819.199900000001
$ws.Range("L132").Value = 11401.5792
$ws.Range("M132").Value = -3289.199900000001
$ws.Range("N132").Value = -16461.5792
$ws.Range("H135").Value = 29775.572
$ws.Range("J135").Value = 29775.572
$ws.Range("L135").Value = 29775.572
$ws.Range("N135").Value = -39915.572
$ws.Range("H136").Value = 1634.7894
$ws.Range("I136").Value = 1031.3914
$ws.Range("J136").Value = 2560
$ws.Range("K136").Value = 3094.1742
$ws.Range("L136").Value = 7680
$ws.Range("M136").Value = -544.1741999999999
$ws.Range("N136").Value = -12780
$ws.Range("H139").Value = 36093.89
$ws.Range("J139").Value = 36093.89
$ws.Range("L139").Value = 36093.89
$ws.Range("N139").Value = -46373.89
